$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Column B (Coin name) updates
$ws.Range("B37").Value = "TheGraph"
$ws.Range("B38").Value = "Dai"
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("B46").Value = "WEMIXToken"

# Column C (Link) updates
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"

# Column D (Price) updates - force text to avoid numeric auto-conversion
Set-TextValue "D2" "67.699.90"
Set-TextValue "D3" "3.538.54"
Set-TextValue "D5" "204.66"
Set-TextValue "D6" "556.67"
Set-TextValue "D7" "3.531.87"
Set-TextValue "D8" "0.610"
Set-TextValue "D10" "64.24"
Set-TextValue "D11" "0.660"
Set-TextValue "D13" "0.0000273"
Set-TextValue "D14" "9.97"
Set-TextValue "D15" "4.085.64"
Set-TextValue "D16" "3.528.28"
Set-TextValue "D18" "18.60"
Set-TextValue "D19" "67.528.93"
Set-TextValue "D20" "11.89"
Set-TextValue "D22" "394.76"
Set-TextValue "D23" "12.17"
Set-TextValue "D24" "4.03"
Set-TextValue "D25" "83.50"
Set-TextValue "D27" "2.86"
Set-TextValue "D28" "12.34"
Set-TextValue "D29" "8.94"
Set-TextValue "D30" "719.24"
Set-TextValue "D31" "31.21"
Set-TextValue "D32" "7.19"
Set-TextValue "D33" "11.80"
Set-TextValue "D34" "64.19"
Set-TextValue "D35" "0.113"
Set-TextValue "D36" "39.00"
Set-TextValue "D37" "0.401"
Set-TextValue "D38" "1.00"
Set-TextValue "D40" "3.05"
Set-TextValue "D41" "0.998"
Set-TextValue "D42" "3.072.01"
Set-TextValue "D43" "0.0₃0689"
Set-TextValue "D44" "2.58"
Set-TextValue "D45" "2.78"
Set-TextValue "D47" "0.0413"
Set-TextValue "D49" "138.40"
Set-TextValue "D50" "8.32"
Set-TextValue "D51" "2.87"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E7").Value = "  -2.79%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +13.12%  "
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("E12").Value = "  -6.64%  "
$ws.Range("E13").Value = "  -7.19%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("E15").Value = "  -3.26%  "
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("E20").Value = "  -5.60%  "
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("E23").Value = "  -7.57%  "
$ws.Range("E24").Value = "  -5.70%  "
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("E30").Value = "  +4.93%  "
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("E32").Value = "  -12.98%  "
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("E36").Value = "  -8.97%  "
$ws.Range("E37").Value = "  -5.97%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  -5.63%  "
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("E43").Value = "  -12.94%  "
$ws.Range("E44").Value = "  -11.35%  "
$ws.Range("E45").Value = "  -8.35%  "
$ws.Range("E46").Value = "  +5.15%  "
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("E50").Value = "  -7.42%  "
$ws.Range("E51").Value = "  -7.60%  "

Write-Host "Applied cryptos list update"